$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.35
$ws.Range("G2").Value = 3.4
$ws.Range("H2").Value = 2.24
$ws.Range("I2").Value = 2.26
$ws.Range("J2").Value = 3.8
$ws.Range("P2").Value = 2.16
$ws.Range("T2").Value = 1.69
$ws.Range("V2").Value = 1.79
$ws.Range("W2").Value = 1.41
$ws.Range("X2").Value = 22
$ws.Range("AA2").Value = 30
$ws.Range("AB2").Value = 14.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AF2").Value = 24
$ws.Range("AG2").Value = 13.5
$ws.Range("AN2").Value = 30

$ws.Range("G3").Value = 3.7
$ws.Range("I3").Value = 2.86
$ws.Range("J3").Value = 3.2
$ws.Range("Q3").Value = 1.71
$ws.Range("T3").Value = 1.73
$ws.Range("V3").Value = 1.54
$ws.Range("W3").Value = 1.37

$ws.Range("F4").Value = 2.6
$ws.Range("J4").Value = 2.78
$ws.Range("K4").Value = 5.3
$ws.Range("Q4").Value = 1.84
$ws.Range("R4").Value = 1.27
$ws.Range("V4").Value = 1.44

$ws.Range("F5").Value = 1.3
$ws.Range("G5").Value = 1.36
$ws.Range("H5").Value = 14.5
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 5.9
$ws.Range("L5").Value = 1.4
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.35
$ws.Range("O5").Value = 1.34
$ws.Range("P5").Value = 1.81
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.3
$ws.Range("S5").Value = 3.65
$ws.Range("T5").Value = 2.6
$ws.Range("U5").Value = 1.49
$ws.Range("W5").Value = 3.75
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 44
$ws.Range("AB5").Value = 7.4
$ws.Range("AC5").Value = 15.5
$ws.Range("AD5").Value = 75
$ws.Range("AG5").Value = 14
$ws.Range("AN5").Value = 8.8

$ws.Range("F6").Value = 1.94
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 3.9
$ws.Range("L6").Value = 1.3
$ws.Range("N6").Value = 4.1
$ws.Range("O6").Value = 1.26
$ws.Range("P6").Value = 2.06
$ws.Range("Q6").Value = 1.78
$ws.Range("R6").Value = 1.41
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 1.71
$ws.Range("U6").Value = 2.2
$ws.Range("X6").Value = 18
$ws.Range("Y6").Value = 20
$ws.Range("AB6").Value = 12
$ws.Range("AE6").Value = 60
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 55

$ws.Range("J7").Value = 4.7
$ws.Range("Q7").Value = 1.46
$ws.Range("S7").Value = 2.14
$ws.Range("U7").Value = 2.16
$ws.Range("X7").Value = 36
$ws.Range("AB7").Value = 13.5
$ws.Range("AF7").Value = 12
$ws.Range("AH7").Value = 23
$ws.Range("AK7").Value = 15
$ws.Range("AL7").Value = 32

$ws.Range("F8").Value = 1.61
$ws.Range("G8").Value = 1.73
$ws.Range("H8").Value = 7.4
$ws.Range("I8").Value = 8.8
$ws.Range("J8").Value = 3.35
$ws.Range("K8").Value = 3.9
$ws.Range("M8").Value = 1.12
$ws.Range("N8").Value = 2.44
$ws.Range("O8").Value = 1.57
$ws.Range("P8").Value = 1.48
$ws.Range("Q8").Value = 2.68
$ws.Range("R8").Value = 1.16
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 2.52
$ws.Range("U8").Value = 1.53
$ws.Range("V8").Value = 1.12
$ws.Range("W8").Value = 2.36
$ws.Range("Y8").Value = 21
$ws.Range("AC8").Value = 11
$ws.Range("AH8").Value = 990
$ws.Range("AJ8").Value = 18
$ws.Range("AK8").Value = 27

$ws.Range("G9").Value = 1.9
$ws.Range("P9").Value = 1.76
$ws.Range("W9").Value = 2.1

$ws.Range("G10").Value = 1.94
$ws.Range("I10").Value = 6.2
$ws.Range("J10").Value = 3.3
$ws.Range("U10").Value = 1.93
$ws.Range("W10").Value = 2.06

$ws.Range("F11").Value = 2.14
$ws.Range("G11").Value = 2.3
$ws.Range("H11").Value = 3.5
$ws.Range("J11").Value = 3.4
$ws.Range("L11").Value = 1.42
$ws.Range("P11").Value = 1.78
$ws.Range("T11").Value = 1.79
$ws.Range("W11").Value = 1.78
$ws.Range("AO11").Value = 60

$ws.Range("F12").Value = 2.34
$ws.Range("K12").Value = 3.5
$ws.Range("P12").Value = 1.74
$ws.Range("Q12").Value = 2.14
$ws.Range("AC12").Value = 9
$ws.Range("AE12").Value = 55
$ws.Range("AG12").Value = 13.5
$ws.Range("AL12").Value = 55

$ws.Range("F13").Value = 1.96
$ws.Range("I13").Value = 4.9
$ws.Range("J13").Value = 3.45
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 1.87
$ws.Range("X13").Value = 990
$ws.Range("AB13").Value = 9.6
$ws.Range("AG13").Value = 1000
